$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Name 5")

# xlCenter / xlBottom (implicit Excel defaults used to drop an explicit
# alignment attribute back to "not set" in the saved xf record) and
# xlLineStyleNone, referenced by value since this host has no Excel
# constants module loaded.
$xlCenter = -4108
$xlBottom = -4107
$xlNone   = -4142

# -------------------------------------------------------------------------
# Row 2: Task-1 -> "Understanding the why of the story(bussiness understanding)"
# -------------------------------------------------------------------------
$ws.Range("C2").Value = "Task-1"
$ws.Range("D2").Value = "Understanding the why of the story(bussiness understanding)"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 3
$ws.Range("G2").Formula = "=E2-F2"

# -------------------------------------------------------------------------
# Row 3: Task-2 -> "Understanding the forward and backward linkages"
# D3 loses its box border (keeps the wrap-text alignment only).
# -------------------------------------------------------------------------
$ws.Range("C3").Value = "Task-2"
$ws.Range("D3").Value = "Understanding the forward and backward linkages"
$ws.Range("D3").Borders.LineStyle = $xlNone
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("G3").Formula = "=(E3-F3)"

# -------------------------------------------------------------------------
# Row 4: Task-3 -> "Understanding the concept of HTML"
# D4 loses its box border (keeps the wrap-text alignment only).
# -------------------------------------------------------------------------
$ws.Range("C4").Value = "Task-3"
$ws.Range("D4").Value = "Understanding the concept of HTML"
$ws.Range("D4").Borders.LineStyle = $xlNone
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Formula = "=(E4-F4)"

# -------------------------------------------------------------------------
# Row 5: Task-4 -> "Understanding the concept of CSS"
# D5 loses its box border (keeps the wrap-text alignment only); F5 becomes
# centered like its neighbours.
# -------------------------------------------------------------------------
$ws.Range("C5").Value = "Task-4"
$ws.Range("D5").Value = "Understanding the concept of CSS"
$ws.Range("D5").Borders.LineStyle = $xlNone
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F5").HorizontalAlignment = $xlCenter
$ws.Range("G5").Formula = "=(E5-F5)"

# -------------------------------------------------------------------------
# Row 6: Task-5 -> "Create structure of login page using html and css "
# E6 drops the vertical-center (back to the plain centered style); F6
# becomes centered.
# -------------------------------------------------------------------------
$ws.Range("C6").Value = "Task-5"
$ws.Range("D6").Value = "Create structure of login page using html and css "
$ws.Range("E6").Value = 4
$ws.Range("E6").VerticalAlignment = $xlBottom
$ws.Range("F6").Value = 3
$ws.Range("F6").HorizontalAlignment = $xlCenter
$ws.Range("G6").Formula = "=E6-F6"

# -------------------------------------------------------------------------
# Row 7: Task-6 -> "Understanding Angular JS" (new task)
# E7 drops the vertical-center; F7 becomes centered (still blank).
# -------------------------------------------------------------------------
$ws.Range("C7").Value = "Task-6"
$ws.Range("D7").Value = "Understanding Angular JS"
$ws.Range("E7").Value = 4
$ws.Range("E7").VerticalAlignment = $xlBottom
$ws.Range("F7").Value = $null
$ws.Range("F7").HorizontalAlignment = $xlCenter
$ws.Range("G7").Formula = "=E7-F7"

# -------------------------------------------------------------------------
# Row 8: Task-7 -> "Integrate structure of login page with angular js functionality" (new task)
# E8 drops the vertical-center. F8 stays plain/blank.
# -------------------------------------------------------------------------
$ws.Range("C8").Value = "Task-7"
$ws.Range("D8").Value = "Integrate structure of login page with angular js functionality"
$ws.Range("E8").Value = 3
$ws.Range("E8").VerticalAlignment = $xlBottom
$ws.Range("F8").Value = $null
$ws.Range("G8").Formula = "=E8-F8"

# B8 (bottom cell of the merged B2:B8 "Story Estimate" block) loses its
# bottom border, matching the plain side-only border used by B3:B7.
$ws.Range("B8").Borders.Item(9).LineStyle = $xlNone

# -------------------------------------------------------------------------
# Row 9: totals row ("Total time"), E9/F9 become centered like the data
# column above them.
# -------------------------------------------------------------------------
$ws.Range("D9").Value = "Total time"
$ws.Range("E9").Formula = "=SUM(E2:E8)"
$ws.Range("E9").HorizontalAlignment = $xlCenter
$ws.Range("F9").Formula = "=SUM(F2:F8)"
$ws.Range("F9").HorizontalAlignment = $xlCenter
$ws.Range("G9").Value = $null

# B2 (Story Estimate total) recalculates automatically via =SUM(E2:E8)

# -------------------------------------------------------------------------
# Selection moves to D9, matching the author's final cursor position.
# -------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D9").Select()
